$d = $word.ActiveDocument

# Update the date line at the top of the document.
$d.Content.Find.Execute("2026-01-06 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2026-01-07 Wednesday", 2) | Out-Null

# Update each division-problem cell in the single table, addressed by
# (row, column) rather than by text, since some old values repeat.
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "30÷5=6, 0"
$t.Cell(1, 2).Range.Text = "98÷2=49, 0"
$t.Cell(1, 3).Range.Text = "57÷8=7, 1"
$t.Cell(1, 4).Range.Text = "31÷8=3, 7"
$t.Cell(1, 5).Range.Text = "91÷2=45, 1"
$t.Cell(5, 1).Range.Text = "42÷5=8, 2"
$t.Cell(5, 2).Range.Text = "64÷7=9, 1"
$t.Cell(5, 3).Range.Text = "59÷5=11, 4"
$t.Cell(5, 4).Range.Text = "57÷4=14, 1"
$t.Cell(5, 5).Range.Text = "59÷5=11, 4"
$t.Cell(9, 1).Range.Text = "81÷2=40, 1"
$t.Cell(9, 2).Range.Text = "92÷5=18, 2"
$t.Cell(9, 3).Range.Text = "86÷7=12, 2"
$t.Cell(9, 4).Range.Text = "78÷3=26, 0"
$t.Cell(9, 5).Range.Text = "18÷4=4, 2"
$t.Cell(13, 1).Range.Text = "88÷4=22, 0"
$t.Cell(13, 2).Range.Text = "99÷4=24, 3"
$t.Cell(13, 3).Range.Text = "17÷7=2, 3"
$t.Cell(13, 4).Range.Text = "69÷9=7, 6"
$t.Cell(13, 5).Range.Text = "35÷2=17, 1"
$t.Cell(17, 1).Range.Text = "35÷4=8, 3"
$t.Cell(17, 2).Range.Text = "54÷8=6, 6"
$t.Cell(17, 3).Range.Text = "60÷8=7, 4"
$t.Cell(17, 4).Range.Text = "98÷3=32, 2"
$t.Cell(17, 5).Range.Text = "93÷6=15, 3"
